# Add a new row (row 92) to the price data sheet, continuing the daily
# price series after the last existing row (row 91, dated 2024-11-01).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 92

$ws.Cells.Item($newRow, 1).Value = "2024-11-02 00:00:00"
$ws.Cells.Item($newRow, 2).Value = 73850
$ws.Cells.Item($newRow, 3).Value = 10340.96
$ws.Cells.Item($newRow, 4).Value = 9151.299999999999
$ws.Cells.Item($newRow, 5).Value = 7.1227
